$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.790.49"
$ws.Range("E2").Value = "  +0.93%  "
$ws.Range("D3").Value = "2.100.09"
$ws.Range("E3").Value = "  +1.03%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'226.88"
$ws.Range("E5").Value = "  -0.79%  "
$ws.Range("D6").Value = "'0.616"
$ws.Range("E6").Value = "  +1.07%  "
$ws.Range("D7").Value = "'62.13"
$ws.Range("E7").Value = "  +3.20%  "
$ws.Range("E9").Value = "  +1.60%  "
$ws.Range("E10").Value = "  +0.87%  "
$ws.Range("E11").Value = "  -0.26%  "
$ws.Range("D12").Value = "'15.82"
$ws.Range("E12").Value = "  +6.53%  "
$ws.Range("D13").Value = "2.411.57"
$ws.Range("E13").Value = "  +1.01%  "
$ws.Range("D14").Value = "'21.98"
$ws.Range("E14").Value = "  -1.62%  "
$ws.Range("E15").Value = "  +2.66%  "
$ws.Range("E16").Value = "  +0.89%  "
$ws.Range("D17").Value = "2.113.56"
$ws.Range("E17").Value = "  +1.65%  "
$ws.Range("D18").Value = "38.804.47"
$ws.Range("E18").Value = "  +1.18%  "
$ws.Range("D19").Value = "'71.45"
$ws.Range("E19").Value = "  +0.46%  "
$ws.Range("E21").Value = "  +1.72%  "
$ws.Range("D22").Value = "'227.23"
$ws.Range("E22").Value = "  +0.99%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").Value = "'2.52"
$ws.Range("E24").Value = "  +5.40%  "
$ws.Range("E25").Value = "  -1.60%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "'9.66"
$ws.Range("E26").Value = "  +2.73%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "'170.55"
$ws.Range("E27").Value = "  +0.38%  "
$ws.Range("D28").Value = "'0.135"
$ws.Range("E28").Value = "  -0.35%  "
$ws.Range("D30").Value = "'19.30"
$ws.Range("E30").Value = "  +1.35%  "
$ws.Range("D31").Value = "'2.53"
$ws.Range("E31").Value = "  +9.18%  "
$ws.Range("E32").Value = "  +0.70%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'4.56"
$ws.Range("E33").Value = "  +1.55%  "
$ws.Range("B34").Value = "THORChain"
$ws.Range("C34").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D34").Value = "'7.19"
$ws.Range("E34").Value = "  +13.55%  "
$ws.Range("D35").Value = "'4.76"
$ws.Range("E35").Value = "  -0.34%  "
$ws.Range("E36").Value = "  +1.45%  "
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("D38").Value = "'3.49"
$ws.Range("E38").Value = "  -1.11%  "
$ws.Range("D39").Value = "'0.999"
$ws.Range("E39").Value = "  -0.19%  "
$ws.Range("D40").Value = "'17.98"
$ws.Range("E40").Value = "  -1.80%  "
$ws.Range("D41").Value = "'0.0228"
$ws.Range("E41").Value = "  +3.36%  "
$ws.Range("D42").Value = "'101.37"
$ws.Range("E42").Value = "  +1.26%  "
$ws.Range("D43").Value = "1.523.83"
$ws.Range("E43").Value = "  -1.01%  "
$ws.Range("E44").Value = "  +7.21%  "
$ws.Range("E45").Value = "  +0.13%  "
$ws.Range("E46").Value = "  +0.94%  "
$ws.Range("D47").Value = "'0.0913"
$ws.Range("E47").Value = "  -0.84%  "
$ws.Range("E48").Value = "  +5.12%  "
$ws.Range("E49").Value = "  +1.91%  "
$ws.Range("E50").Value = "  -0.79%  "
$ws.Range("D51").Value = "2.298.45"
